# Daily attendance processing - 2025-10-05 06:46:00
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a text value that LOOKS like a percentage (e.g. "27.5%")
# without letting Excel auto-convert it into a numeric percent cell.
# We briefly mark the cell as Text (@), assign the literal string, then
# restore the cell's original number format/style by pasting the format
# from a neighbouring cell that already carries the untouched style, so
# the cell's style index is left exactly as it was.
function Set-TextValue {
    param($cellRef, $text, $formatDonorRef)
    $target = $ws.Range($cellRef)
    $target.NumberFormat = "@"
    $target.Value = $text
    $ws.Range($formatDonorRef).Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null
}

# Row 6 - Recorded Sessions
$ws.Range("L6").Value = 42

# Row 7 - PARASITOLOGY Recorded By order swap
$ws.Range("G7").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 8 - Pending Sessions
$ws.Range("L8").Value = 103

# Row 9 - Coverage %
Set-TextValue "L9" "27.5%" "L10"

# Row 10 - Average Attendance %
Set-TextValue "L10" "48.0%" "L9"

# Row 14 - Recorded By order swap
$ws.Range("G14").Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"

# Row 17 - Recorded By order swap
$ws.Range("G17").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# Row 18 - Group A4 statistics
$ws.Range("O18").Value = 7
$ws.Range("Q18").Value = 10
Set-TextValue "R18" "41.2%" "S18"
Set-TextValue "S18" "47.7%" "R18"

# Row 24 - Recorded By order swap
$ws.Range("G24").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 31 - Recorded By order swap
$ws.Range("G31").Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"

# Row 34 - Recorded By order swap
$ws.Range("G34").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# Row 35 - Recorded By order swap
$ws.Range("G35").Value = "Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"

# Row 45 - Recorded By order swap
$ws.Range("G45").Value = "backup@backdoor.com, System, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 51 - Recorded By order swap
$ws.Range("G51").Value = "abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

# Row 62 - Recorded By order swap
$ws.Range("G62").Value = "backup@backdoor.com, System, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 68 - Recorded By order swap
$ws.Range("G68").Value = "abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

# Row 69 - PHYSIOLOGY Y2-A4 session 2 now recorded (copy the "Recorded" look
# from row 68 first so fill/style match exactly, then set its own values)
$ws.Range("A68:I68").Copy() | Out-Null
$ws.Range("A69:I69").PasteSpecial(-4122) | Out-Null
$ws.Range("G69").Value = "yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("H69").Value = "23/225"
$ws.Range("I69").Value = "Recorded"

# Row 72 - Recorded By order swap
$ws.Range("G72").Value = "Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"

# Row 83 - Recorded By order swap
$ws.Range("G83").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"

# Row 85 - Recorded By order swap
$ws.Range("G85").Value = "wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# Row 98 - Recorded By order swap
$ws.Range("G98").Value = "afaf.abdallah@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, user@user.com, nourhanmohamed@med.asu.edu.eg"

# Row 99 - Recorded By order swap
$ws.Range("G99").Value = "Walaa.h.ghanima@med.asu.edu.eg, user@user.com"

# Row 102 - Recorded By order swap
$ws.Range("G102").Value = "wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# Row 109 - Recorded By order swap
$ws.Range("G109").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 119 - Recorded By order swap
$ws.Range("G119").Value = "shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# Row 126 - Recorded By order swap
$ws.Range("G126").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 136 - Recorded By order swap
$ws.Range("G136").Value = "shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# Row 143 - Recorded By order swap
$ws.Range("G143").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 149 - Recorded By order swap
$ws.Range("G149").Value = "Walaa.h.ghanima@med.asu.edu.eg, user@user.com"

# Row 150 - Recorded By order swap
$ws.Range("G150").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
